# Add 8 new "Summary_*" data-source sheets at the end of the workbook,
# continuing the existing repeating 4-sheet pattern:
#   1) TestData breakdown, 6 rows   (Name/Age/City/Score/Comments vs TestData!A..E)
#   2) header-only sheet, 1 row
#   3) TestData breakdown, 7 rows   (adds a 6th, unlabeled row vs TestData!F)
#   4) Summary_1 breakdown, 6 rows  (vs Summary_1!A..E)

$wb = $excel.ActiveWorkbook

$headers = @("Column Heading", "Count", "Total", "Percentage")

function Fill-Header($ws) {
    $ws.Range("A1").Value = $headers[0]
    $ws.Range("B1").Value = $headers[1]
    $ws.Range("C1").Value = $headers[2]
    $ws.Range("D1").Value = $headers[3]
}

function Fill-DCol($ws, $row) {
    $ws.Cells.Item($row, 4).Formula = "=B$row/C$row"
    $ws.Cells.Item($row, 4).NumberFormat = "0.00%"
}

# Pattern A: breakdown of TestData, 6 rows
function Fill-TestData6($ws) {
    Fill-Header $ws

    $labels = @("Name", "Age", "City", "Score", "Comments")
    $cols   = @("A", "B", "C", "D", "E")
    for ($i = 0; $i -lt 5; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $labels[$i]
        $ws.Cells.Item($row, 2).Formula = "=COUNTA(TestData!$($cols[$i]):$($cols[$i]))-1"
        $ws.Cells.Item($row, 3).Formula = "=COUNTA(TestData!A:A)-1"
        Fill-DCol $ws $row
    }
}

# Pattern B: header row only
function Fill-EmptyHeader($ws) {
    Fill-Header $ws
}

# Pattern C: breakdown of TestData, 7 rows (adds unlabeled F-column row)
function Fill-TestData7($ws) {
    Fill-TestData6 $ws

    $ws.Cells.Item(7, 2).Formula = "=COUNTA(TestData!F:F)-1"
    $ws.Cells.Item(7, 3).Formula = "=COUNTA(TestData!A:A)-1"
    Fill-DCol $ws 7
}

# Pattern D: breakdown of Summary_1, 6 rows
function Fill-Summary1_6($ws) {
    Fill-Header $ws

    $labels = @("Column Heading", "Count", "Total", "Percentage")
    $cols   = @("A", "B", "C", "D")
    for ($i = 0; $i -lt 4; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $labels[$i]
        $ws.Cells.Item($row, 2).Formula = "=COUNTA(Summary_1!$($cols[$i]):$($cols[$i]))-1"
        $ws.Cells.Item($row, 3).Formula = "=COUNTA(Summary_1!A:A)-1"
        Fill-DCol $ws $row
    }

    # row 6 has no label in column A
    $ws.Cells.Item(6, 2).Formula = "=COUNTA(Summary_1!E:E)-1"
    $ws.Cells.Item(6, 3).Formula = "=COUNTA(Summary_1!A:A)-1"
    Fill-DCol $ws 6
}

$pattern = @("TestData6", "EmptyHeader", "TestData7", "Summary1_6")

for ($n = 50; $n -le 57; $n++) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = "Summary_$n"

    $kind = $pattern[($n - 50) % 4]
    switch ($kind) {
        "TestData6"    { Fill-TestData6 $newSheet }
        "EmptyHeader"  { Fill-EmptyHeader $newSheet }
        "TestData7"    { Fill-TestData7 $newSheet }
        "Summary1_6"   { Fill-Summary1_6 $newSheet }
    }
}
